# Inserts a new client row ("CORPORACION TOAQUIZAJR CORPOTOAQUIJR S.A.S.")
# at row 9 (alphabetically before "FREVIUNO CIA. LTDA.") in both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, pushing the existing rows
# down by one, and refreshes the trailing summary row so it reflects the
# new row count (28 -> 29 clients).

$wb = $excel.ActiveWorkbook
$advisor = "ALMEIDA CUATIN JHONATHANN CARLOS"
$newClient = "CORPORACION TOAQUIZAJR CORPOTOAQUIJR S.A.S."

# --- Sheet "VENTAS POR GRUPO" (columns A:R, data rows 2-29, summary row 30) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(9).Insert()
$ws1.Cells.Item(9, 1).Value = $advisor
$ws1.Cells.Item(9, 2).Value = $newClient
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(9, $c).Value = 0
}
# Summary row moved from 30 -> 31; update the "X de 28" counters to "X de 29"
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(31, $c)
    $cell.Value = $cell.Value().Replace("de 28", "de 29")
}

# --- Sheet "VENTA MENSUAL" (columns A:G, data rows 2-29, summary row 30) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(9).Insert()
$ws2.Cells.Item(9, 1).Value = $advisor
$ws2.Cells.Item(9, 2).Value = $newClient
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(9, $c).Value = 0
}
# Summary row (plain numeric totals) moved from 30 -> 31; totals are unchanged
# since the inserted row contributes only zeros.
